$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "First name"
$ws.Range("B1").Value = "Last name"
$ws.Range("C1").Value = "Email"

$ws.Range("C2").Select()
